$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -3.380387134729135
$ws.Range("C2").Value = 2.152106832695933
$ws.Range("D2").Value = 7.434863256980662

$ws.Range("B3").Value = 0.2458182668426012
$ws.Range("C3").Value = -0.175034069307578
$ws.Range("D3").Value = -1.006965064779253

$ws.Range("B4").Value = 5.628575150153137
$ws.Range("C4").Value = 0.6539026922407265
$ws.Range("D4").Value = 6.267597251366408

$ws.Range("B5").Value = 6.17726467992108
$ws.Range("C5").Value = -6.90125146503876
$ws.Range("D5").Value = 9.755157674477211

$ws.Range("B6").Value = -1.454355621792969
$ws.Range("C6").Value = -6.17799605785695
$ws.Range("D6").Value = 8.002872036043618

$ws.Range("B7").Value = -0.2862285124544894
$ws.Range("C7").Value = -4.971795880184382
$ws.Range("D7").Value = 2.79976609910868

$ws.Range("B8").Value = -0.9025313378329569
$ws.Range("C8").Value = -4.255046939480378
$ws.Range("D8").Value = 0.1133182134140931

$ws.Range("B9").Value = 4.667723483238428
$ws.Range("C9").Value = -1.176744080281222
$ws.Range("D9").Value = 11.23837728813852

$ws.Range("B10").Value = -10.66418995593149
$ws.Range("C10").Value = -5.569625683566737
$ws.Range("D10").Value = -5.930111789747738

$ws.Range("B11").Value = -6.409519563932187
$ws.Range("C11").Value = 8.992974645904539
$ws.Range("D11").Value = -7.356937051210577

$ws.Range("B12").Value = 0.07754669076678322
$ws.Range("C12").Value = 7.485228308347747
$ws.Range("D12").Value = -10.82673191893189

$ws.Range("B13").Value = -2.471908895400521
$ws.Range("C13").Value = 2.942865220470381
$ws.Range("D13").Value = -3.253364144377369
